# Updated cryptos list on Thu Feb 29 13:35:50 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for the cryptos
# table, and swaps the Stellar/Monero rows (42/43) to reflect their new
# rank order (with their Coin/Link/Price/Volume values moving together).
#
# D-column values are written with a leading apostrophe so Excel keeps
# them as literal text (e.g. "407.81", "1.99", "3.60") instead of
# re-interpreting them as numbers and dropping trailing zeros / changing
# precision. E-column values already contain non-numeric padding/percent
# text so they naturally stay text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.431.92"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").Value = "'3.465.01"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'407.81"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "'131.53"
$ws.Range("E6").Value = "  +17.35%  "
$ws.Range("D7").Value = "'3.464.31"
$ws.Range("E7").Value = "  +3.86%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'0.690"
$ws.Range("E10").Value = "  +8.91%  "
$ws.Range("D11").Value = "'0.128"
$ws.Range("E11").Value = "  +29.87%  "
$ws.Range("D12").Value = "'42.81"
$ws.Range("E12").Value = "  +6.78%  "
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "'4.008.24"
$ws.Range("E14").Value = "  +3.51%  "
$ws.Range("D15").Value = "'8.72"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").Value = "'20.07"
$ws.Range("E16").Value = "  +3.58%  "
$ws.Range("D17").Value = "'3.447.67"
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("D18").Value = "'62.390.85"
$ws.Range("E18").Value = "  +3.90%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "'10.87"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "'0.0000137"
$ws.Range("E21").Value = "  +24.15%  "
$ws.Range("D23").Value = "'82.71"
$ws.Range("E23").Value = "  +9.29%  "
$ws.Range("D24").Value = "'13.12"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'308.68"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "'30.31"
$ws.Range("E27").Value = "  +5.76%  "
$ws.Range("D28").Value = "'8.27"
$ws.Range("E28").Value = "  +4.85%  "
$ws.Range("D29").Value = "'7.75"
$ws.Range("E29").Value = "  +3.89%  "
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("D31").Value = "'4.37"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("E32").Value = "  +4.25%  "
$ws.Range("D33").Value = "'2.66"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'11.92"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").Value = "'43.01"
$ws.Range("E35").Value = "  +8.10%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'0.0492"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("D38").Value = "'52.63"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").Value = "'3.60"
$ws.Range("E39").Value = "  +5.76%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'2.99"
$ws.Range("E41").Value = "  -6.26%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.126"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'137.87"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  +3.62%  "
$ws.Range("D45").Value = "'17.54"
$ws.Range("E45").Value = "  +3.73%  "
$ws.Range("D46").Value = "'3.97"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "'0.286"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'22.38"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "'2.201.40"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "'3.806.27"
$ws.Range("E51").Value = "  +3.94%  "
